$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename labels to reflect the AutoMOF_6 series wording
$ws.Range("F4").Value = " Zn6H_2"
$ws.Range("F6").Value = " MeOH_2"
$ws.Range("E3").Value = " mass of Zn6H_2"
$ws.Range("E5").Value = " mass of MeOH_2"

# Leave the final selection on E5, matching the edited cell that was last touched
$ws.Range("E5").Select()
